# Fix the FN(%)/FP(%)/VN(%)/VP(%) summary columns on the "general_summary"
# sheet: they were computed with the wrong denominator (total samples)
# instead of the matching pair (e.g. FN(%) = FN / (FN + VN)). A leading
# apostrophe forces the corrected value to be stored as text, matching the
# original (numeric-looking) text cells in these columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general_summary")

$ws.Range("H2").Value = "'100.0"
$ws.Range("J2").Value = "'100.0"

$ws.Range("H3").Value = "'40.0"
$ws.Range("J3").Value = "'100.0"
$ws.Range("L3").Value = "'60.0"

$ws.Range("H4").Value = "'60.7"
$ws.Range("J4").Value = "'71.3"
$ws.Range("L4").Value = "'39.3"
$ws.Range("N4").Value = "'28.7"

$ws.Range("H5").Value = "'62.3"
$ws.Range("J5").Value = "'66.7"
$ws.Range("L5").Value = "'37.7"
$ws.Range("N5").Value = "'33.3"

$ws.Range("H6").Value = "'56.1"
$ws.Range("J6").Value = "'51.9"
$ws.Range("L6").Value = "'43.9"
$ws.Range("N6").Value = "'48.1"

$ws.Range("H7").Value = "'61.5"
$ws.Range("J7").Value = "'87.2"
$ws.Range("L7").Value = "'38.5"
$ws.Range("N7").Value = "'12.8"
